$wb = $excel.ActiveWorkbook

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2333.3333
$ws.Range("I40").Value = 2333.3333
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2333.3333
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2158.3333
$ws.Range("N40").ClearContents()

# ALC row 42
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 989.3333
$ws.Range("I42").Value = 21.666666
$ws.Range("J42").Value = 1473.1666
$ws.Range("K42").Value = 64.99999800000001
$ws.Range("L42").Value = 4419.4998
$ws.Range("M42").Value = 165.000002
$ws.Range("N42").Value = -4879.4998

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 9000
$ws.Range("J51").Value = 7000
$ws.Range("L51").Value = 7000
$ws.Range("N51").Value = -7968

# ALC row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 49281.43
$ws.Range("J87").Value = 49281.43
$ws.Range("L87").Value = 49281.43
$ws.Range("N87").Value = -51777.43

# ALC row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 49281.43
$ws.Range("J90").Value = 49281.43
$ws.Range("L90").Value = 147844.29
$ws.Range("N90").Value = -160324.29

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 557
$ws.Range("I111").Value = 557
$ws.Range("K111").Value = 1671
$ws.Range("M111").Value = 1396

# ALC row 115
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 375.66666
$ws.Range("I115").Value = 375.66666
$ws.Range("K115").Value = 1126.99998
$ws.Range("M115").Value = 440.0000199999999

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6226.5
$ws.Range("I116").Value = 5998
$ws.Range("J116").Value = 6455
$ws.Range("K116").Value = 5998
$ws.Range("L116").Value = 6455
$ws.Range("M116").Value = -2556
$ws.Range("N116").Value = -13339

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2637.9375
$ws.Range("I132").Value = 2637.9375
$ws.Range("K132").Value = 7913.8125
$ws.Range("M132").Value = -5383.8125

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1997.1904
$ws.Range("I138").Value = 1194.3
$ws.Range("J138").Value = 2727.0908
$ws.Range("K138").Value = 3582.9
$ws.Range("L138").Value = 8181.2724
$ws.Range("M138").Value = 1557.1
$ws.Range("N138").Value = -18461.2724

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1584.2307
$ws.Range("I141").Value = 1584.2307
$ws.Range("K141").Value = 4752.6921
$ws.Range("M141").Value = 427.3078999999998

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1404.909
$ws.Range("I32").Value = 1225.3442
$ws.Range("J32").Value = 3595.6
$ws.Range("K32").Value = 1225.3442
$ws.Range("L32").Value = 3595.6
$ws.Range("M32").Value = -938.3442
$ws.Range("N32").Value = -4169.6

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1092.7273
$ws.Range("I74").Value = 1007.7778
$ws.Range("K74").Value = 1007.7778
$ws.Range("M74").Value = -133.7778

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1092.7273
$ws.Range("I77").Value = 1007.7778
$ws.Range("K77").Value = 5038.889
$ws.Range("M77").Value = -670.8890000000001

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4315.1665
$ws.Range("I20").Value = 3778.6
$ws.Range("K20").Value = 3778.6
$ws.Range("M20").Value = -3531.6

# CRP row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 21181.5
$ws.Range("J51").Value = 24999.8
$ws.Range("L51").Value = 24999.8
$ws.Range("N51").Value = -26471.8

# CRP row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 21181.5
$ws.Range("J61").Value = 24999.8
$ws.Range("L61").Value = 24999.8
$ws.Range("N61").Value = -25695.8

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 698.4
$ws.Range("I107").Value = 298.5
$ws.Range("J107").Value = 965
$ws.Range("K107").Value = 298.5
$ws.Range("L107").Value = 965
$ws.Range("M107").Value = 1621.5
$ws.Range("N107").Value = -4805

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2539.5
$ws.Range("I5").Value = 3164
$ws.Range("J5").Value = 1915
$ws.Range("K5").Value = 9492
$ws.Range("L5").Value = 5745
$ws.Range("M5").Value = -9380
$ws.Range("N5").Value = -5969

# CUL row 64
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

# CUL row 67
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1055.4
$ws.Range("J121").Value = 1015.75
$ws.Range("L121").Value = 3047.25
$ws.Range("N121").Value = -5667.25

# CUL row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 127121.875
$ws.Range("I134").Value = 166990.33
$ws.Range("J134").Value = 7516.5
$ws.Range("K134").Value = 500970.99
$ws.Range("L134").Value = 22549.5
$ws.Range("M134").Value = -495900.99
$ws.Range("N134").Value = -32689.5

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2539.5
$ws.Range("I135").Value = 3164
$ws.Range("J135").Value = 1915
$ws.Range("K135").Value = 28476
$ws.Range("L135").Value = 17235
$ws.Range("M135").Value = -25941
$ws.Range("N135").Value = -22305

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 333333340
$ws.Range("I70").Value = 333333340
$ws.Range("K70").Value = 333333340
$ws.Range("M70").Value = -333333070

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 333333340
$ws.Range("I73").Value = 333333340
$ws.Range("K73").Value = 333333340
$ws.Range("M73").Value = -333332404

# GSM row 96
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 40000
$ws.Range("J96").Value = 40000
$ws.Range("L96").Value = 40000
$ws.Range("N96").Value = -45492

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 12287.125
$ws.Range("I126").Value = 12287.125
$ws.Range("K126").Value = 36861.375
$ws.Range("M126").Value = -34391.375

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2364.8235
$ws.Range("I132").Value = 2393.7334
$ws.Range("K132").Value = 7181.2002
$ws.Range("M132").Value = -4651.2002

# LTW row 120
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H120").Value = 59998
$ws.Range("J120").Value = 59998
$ws.Range("L120").Value = 59998
$ws.Range("N120").Value = -69674

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 14057.958
$ws.Range("I132").Value = 9820.786
$ws.Range("K132").Value = 29462.358
$ws.Range("M132").Value = -26932.358

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2821.3635
$ws.Range("I136").Value = 2411.7856
$ws.Range("K136").Value = 7235.3568
$ws.Range("M136").Value = -4685.3568

# WVR row 95
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 19500
$ws.Range("J95").Value = 19500
$ws.Range("L95").Value = 19500
$ws.Range("N95").Value = -24992

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3488
$ws.Range("I132").Value = 3488
$ws.Range("K132").Value = 10464
$ws.Range("M132").Value = -7934

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4177.1816
$ws.Range("I136").Value = 4281
$ws.Range("J136").Value = 3995.5
$ws.Range("K136").Value = 12843
$ws.Range("L136").Value = 11986.5
$ws.Range("M136").Value = -10293
$ws.Range("N136").Value = -17086.5
